$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '98.767.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.26%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.346.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.78%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '651.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.62%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.55'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +11.60%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.466'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +18.45%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.09'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +24.48%  '

$ws.Range('E10').Value = '  +0.02%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.341.76'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.87%  '

$ws.Range('E12').Value = '  +5.27%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '43.48'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +20.50%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000269'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.44%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '98.635.72'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.11%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.982.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.07%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.09%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.345.27'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.51%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +23.38%  '

$ws.Range('E20').Value = '  +10.87%  '

$ws.Range('B21').Value = 'SuiNetwork'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.61'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.29%  '

$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '536.97'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.72%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.34'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.88%  '

$ws.Range('E24').Value = '  +1.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.440'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +56.35%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '102.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +15.86%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +10.89%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.97%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.526.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.00%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.151'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +13.65%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.03'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +14.90%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.191'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.20%  '

$ws.Range('E34').Value = '  +0.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '29.44'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.73%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.537'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +15.89%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.86'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.89%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.10'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.33%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.156'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.07%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '521.77'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.78%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '24.71'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.65%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.18%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.77'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.43%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0425'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +30.06%  '

$ws.Range('E45').Value = '  +3.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.823'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.23%  '

$ws.Range('E47').Value = '  +0.00%  '

$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.89'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +20.03%  '

$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.46%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.08'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.74%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '164.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.42%  '
